$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 168
$ws.Range("B1").Value = 214.1999999999989
$ws.Range("C1").Value = 268.5999999999986
$ws.Range("A2").Value = 168
